$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6133912801742554
$ws.Range("B1").Value = 1.355645537376404
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.35623025894165
$ws.Range("E1").Value = 1.364726424217224
